$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.945.91'
$ws.Range("E2").Value = '  +4.45%  '
$ws.Range("D3").Value = '3.257.26'
$ws.Range("E3").Value = '  +2.57%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '395.99'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.98%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '108.86'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.41%  '
$ws.Range("E7").Value = '  +7.15%  '
$ws.Range("D8").Value = '3.255.61'
$ws.Range("E8").Value = '  +2.73%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.628'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +2.08%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '39.27'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.28%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0988'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +11.84%  '
$ws.Range("E13").Value = '  +2.17%  '
$ws.Range("D14").Value = '3.776.74'
$ws.Range("E14").Value = '  +2.71%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.39'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +4.55%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '19.20'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.63%  '
$ws.Range("D17").Value = '3.258.33'
$ws.Range("E17").Value = '  +2.64%  '
$ws.Range("E18").Value = '  -2.74%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.76'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.96%  '
$ws.Range("D20").Value = '56.878.68'
$ws.Range("E20").Value = '  +4.55%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '3.34'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +1.74%  '
$ws.Range("E22").Value = '  +9.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.05'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +1.28%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '296.18'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +7.62%  '
$ws.Range("E25").Value = '  +3.03%  '
$ws.Range("E26").Value = '  -2.28%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '28.18'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.41%  '
$ws.Range("E28").Value = '  +1.04%  '
$ws.Range("E29").Value = '  -3.90%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.25'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -4.25%  '
$ws.Range("E31").Value = '  -0.63%  '
$ws.Range("E32").Value = '  +0.08%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.26'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +2.11%  '
$ws.Range("E34").Value = '  -2.78%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '39.89'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +9.07%  '
$ws.Range("E36").Value = '  -3.39%  '
$ws.Range("E37").Value = '  +1.49%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '51.49'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.15%  '
$ws.Range("E39").Value = '  -0.04%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.49'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -3.51%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.96'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +2.60%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '139.52'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +6.29%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.123'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +4.36%  '
$ws.Range("B44").Value = 'ARBITRUM'
$ws.Range("C44").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.90'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -1.94%  '
$ws.Range("B45").Value = 'NEARProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.99'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -2.52%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '17.18'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.15%  '
$ws.Range("E47").Value = '  -3.33%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.25'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.76%  '
$ws.Range("E49").Value = '  +3.43%  '
$ws.Range("D50").Value = '2.165.83'
$ws.Range("E50").Value = '  +3.49%  '
$ws.Range("E51").Value = '  -6.36%  '
